# Updated cryptos list values (prices and 1h volume %) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'66.308.77"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Formula = "'3.061.01"
$ws.Range("E3").Value = "  +2.44%  "
$ws.Range("D4").Formula = "'0.999"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Formula = "'578.76"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").Formula = "'168.30"
$ws.Range("E6").Value = "  +4.23%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Formula = "'3.057.56"
$ws.Range("E8").Value = "  +2.43%  "
$ws.Range("D9").Formula = "'0.523"
$ws.Range("E9").Value = "  +1.34%  "
$ws.Range("D10").Formula = "'6.70"
$ws.Range("E10").Value = "  +3.17%  "
$ws.Range("D11").Formula = "'0.153"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").Formula = "'0.486"
$ws.Range("E12").Value = "  +7.34%  "
$ws.Range("D13").Formula = "'0.0000249"
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").Formula = "'36.87"
$ws.Range("E14").Value = "  +7.71%  "
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("D16").Formula = "'66.329.98"
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").Formula = "'3.553.86"
$ws.Range("E17").Value = "  +2.02%  "
$ws.Range("D18").Formula = "'7.24"
$ws.Range("E18").Value = "  +5.79%  "
$ws.Range("D19").Formula = "'16.70"
$ws.Range("E19").Value = "  +21.42%  "
$ws.Range("D20").Formula = "'3.048.65"
$ws.Range("E20").Value = "  +2.16%  "
$ws.Range("D21").Formula = "'466.66"
$ws.Range("E21").Value = "  +3.71%  "
$ws.Range("D22").Formula = "'0.712"
$ws.Range("E22").Value = "  +4.86%  "
$ws.Range("D23").Formula = "'7.45"
$ws.Range("E23").Value = "  +2.64%  "
$ws.Range("D24").Formula = "'83.22"
$ws.Range("E24").Value = "  +1.59%  "
$ws.Range("D25").Formula = "'12.88"
$ws.Range("E25").Value = "  +5.95%  "
$ws.Range("D26").Formula = "'2.29"
$ws.Range("E26").Value = "  +1.97%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Formula = "'10.05"
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Formula = "'1.00"
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").Formula = "'8.18"
$ws.Range("E29").Value = "  +1.13%  "
$ws.Range("E30").Value = "  +2.30%  "
$ws.Range("D31").Formula = "'2.67"
$ws.Range("E31").Value = "  +2.45%  "
$ws.Range("E32").Value = "  -1.05%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Formula = "'28.35"
$ws.Range("E33").Value = "  +4.95%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Formula = "'0.116"
$ws.Range("E34").Value = "  +6.24%  "
$ws.Range("D35").Formula = "'1.00"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").Formula = "'0.998"
$ws.Range("E36").Value = "  +1.44%  "
$ws.Range("D37").Formula = "'5.88"
$ws.Range("E37").Value = "  +2.06%  "
$ws.Range("D38").Formula = "'48.22"
$ws.Range("E38").Value = "  +9.88%  "
$ws.Range("D39").Formula = "'0.321"
$ws.Range("E39").Value = "  +6.65%  "
$ws.Range("D40").Formula = "'49.76"
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("E41").Value = "  -1.68%  "
$ws.Range("D42").Formula = "'0.122"
$ws.Range("E42").Value = "  +0.94%  "
$ws.Range("D43").Formula = "'8.68"
$ws.Range("E43").Value = "  +3.70%  "
$ws.Range("D44").Formula = "'2.84"
$ws.Range("E44").Value = "  -1.19%  "
$ws.Range("D45").Formula = "'0.0361"
$ws.Range("E45").Value = "  +1.19%  "
$ws.Range("D46").Formula = "'382.77"
$ws.Range("E46").Value = "  -2.97%  "
$ws.Range("D47").Formula = "'2.744.86"
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("D48").Formula = "'133.84"
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").Formula = "'24.57"
$ws.Range("E50").Value = "  +5.20%  "
$ws.Range("D51").Formula = "'2.23"
$ws.Range("E51").Value = "  +4.53%  "
